$d = $word.ActiveDocument

function Get-ParagraphRangeByText($doc, $searchText) {
    $r = $doc.Content
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Could not find text: " + $searchText)
    }
    $para = $r.Paragraphs(1)
    return $para.Range
}

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$colorRPr = '<w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr>'

# --- Item 4: "Registrierungsbestätigung via Email" -> split "Email" off with gramStart/gramEnd proofErr ---
$r1 = Get-ParagraphRangeByText $d "Registrierungsbestätigung via Email"
$xml1 = '<w:p ' + $wns + '><w:r><w:t>4</w:t></w:r><w:r><w:tab/></w:r>' + `
    '<w:r>' + $colorRPr + '<w:t xml:space="preserve">Registrierungsbestätigung via </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $colorRPr + '<w:t>Email</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$r1.InsertXML($xml1)

# --- Item 10: "Artikel hinzufügen und ändern der Menge im Warenkorb" -> split with _GoBack bookmark ---
$r2 = Get-ParagraphRangeByText $d "Artikel hinzufügen und ändern der Menge im Warenkorb"
$xml2 = '<w:p ' + $wns + '><w:pPr><w:ind w:left="708" w:hanging="708"/><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr>' + `
    '<w:r><w:t>10</w:t></w:r>' + `
    '<w:r>' + $colorRPr + '<w:tab/></w:r>' + `
    '<w:r>' + $colorRPr + '<w:t>Artikel hinzufügen und ändern der Menge i</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r>' + $colorRPr + '<w:t>m Warenkorb</w:t></w:r>' + `
    '</w:p>'
$r2.InsertXML($xml2)

# --- Item 12: shorten text, drop highlight/extra sentences/symbol/bookmark, split "Email" with proofErr ---
$r3 = Get-ParagraphRangeByText $d "Bestätigung via Email (Funktion) mit Artikelmenge + Artikelname + Gesamtsumme + Versandkosten"
$xml3 = '<w:p ' + $wns + '><w:pPr><w:ind w:left="708" w:hanging="708"/></w:pPr>' + `
    '<w:r><w:t>12</w:t></w:r>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:r>' + $colorRPr + '<w:t xml:space="preserve">Bestätigung via </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $colorRPr + '<w:t>Email</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $colorRPr + '<w:t xml:space="preserve"> (Funktion) mit Artikelmenge + Artikelname + Gesamtsumme + Versandkosten</w:t></w:r>' + `
    '</w:p>'
$r3.InsertXML($xml3)

# --- Item 15: "Auf der Startseite ein Karussell - Example: (...)" -> split "Example" with spellStart/spellEnd proofErr ---
$r4 = Get-ParagraphRangeByText $d "Auf der Startseite ein Karussell - Example: (https://getbootstrap.com/docs/4.3/components/carousel/)"
$xml4 = '<w:p ' + $wns + '><w:pPr><w:ind w:left="708" w:hanging="708"/><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr>' + `
    '<w:r><w:t>15</w:t></w:r>' + `
    '<w:r>' + $colorRPr + '<w:tab/><w:t xml:space="preserve">Auf der Startseite ein Karussell - </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $colorRPr + '<w:t>Example</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $colorRPr + '<w:t>: (https://getbootstrap.com/docs/4.3/components/carousel/)</w:t></w:r>' + `
    '</w:p>'
$r4.InsertXML($xml4)

Write-Output "Done"
